$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Coco & Eve Super Hydration Kit
# (shared-string creation order: Name, Link, Image, Price)
$ws.Range("A11").Value = "Coco & Eve`nSUPER HYDRATION KIT"
$ws.Range("C11").Value = "https://fr.zalando.ch/coco-and-eve-super-hydration-kit250ml-shampoo-250ml-conditioner-set-pour-les-cheveux-c1o34h002-s11.html?_rfl=de"
$ws.Range("B11").Value = "https://img01.ztat.net/article/spp-media-p1/9d02fd67425c4703b5bb9a29dfdb234d/da3a0877a08249b7a44ee5bb4d24335f.jpg?imwidth=1800&filter=packshot"
$ws.Range("D11").Value = "37 CHF"
$ws.Range("A11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 30

# Row 12 - MAC Matte Lipstick
# (shared-string creation order: Name, Image, Link, Price)
$ws.Range("A12").Value = "MAC`nMATTE LIPSTICK - culoare mehr"
$ws.Range("B12").Value = "https://img01.ztat.net/article/spp-media-p1/4c4155d867594bb7969d7b3c5d8b5ce9/4815417c961e4ab4848f45d7a7a4bf2a.jpg?imwidth=1800&filter=packshot"
$ws.Range("C12").Value = "https://fr.zalando.ch/mac-matte-lipstick-richard-quinn-exclusive-edition-rouge-a-levres-mehr-m3t31e094-j15.html?_rfl=de"
$ws.Range("D12").Value = "25 CHF"
$ws.Range("A12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 45

# Update the active selection to reflect where the user ended up (one row below the new data)
$ws.Range("D13").Select() | Out-Null
